$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw input values (dependent formulas recalculate automatically)
$ws.Range("K16").Value = 260
$ws.Range("L16").Value = 3153
$ws.Range("E17").Value = 49
$ws.Range("K17").Value = 394

# Update the selection to match the new active cell / selection range
$ws.Range("K18").Select()
